# Weekly update: insert a new price record for Espinaca at
# "Mercado Mayorista Lo Valledor de Santiago" right before the existing
# row 435, shifting the remaining historical rows (435-469) down by one
# (to 436-470). This matches the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 435; Excel shifts rows 435:469 down to 436:470
# and extends the used range (dimension becomes A1:R470) automatically.
$ws.Rows.Item(435).Insert()

# Populate the newly inserted row 435 with the new weekly record.
$ws.Cells.Item(435, 1).Value = 6
$ws.Cells.Item(435, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(435, 3).Value = 'Metropolitana'
$ws.Cells.Item(435, 4).Value = 44578
$ws.Cells.Item(435, 5).Value = 13
$ws.Cells.Item(435, 6).Value = 100112012
$ws.Cells.Item(435, 7).Value = 'Espinaca'
$ws.Cells.Item(435, 8).Value = 'Sin especificar'
$ws.Cells.Item(435, 9).Value = 'Primera'
$ws.Cells.Item(435, 10).Value = 330
$ws.Cells.Item(435, 11).Value = 9500
$ws.Cells.Item(435, 12).Value = 10000
$ws.Cells.Item(435, 13).Value = 9697
$ws.Cells.Item(435, 14).Value = '$/cuna 10 kilos'
$ws.Cells.Item(435, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(435, 16).Value = 970
$ws.Cells.Item(435, 17).Value = 10
$ws.Cells.Item(435, 18).Value = 'Hortaliza'

# Keep the date column formatted the same as the rest of column D
# (numFmt used throughout the sheet for dates).
$ws.Cells.Item(435, 4).NumberFormat = $ws.Cells.Item(436, 4).NumberFormat()
